{"js": "// Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line,\n// the \"\u00a9 2020 . Contact: ...\" credits line, and the blank paragraph that\n// separates them from the preceding \"Requisitos\" content, mirroring the\n// upstream Jekyll/Github Pages rebuild that dropped that boilerplate.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\nconst VER_TEXT = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst COPYRIGHT_MARKER = \"Contact: luizeleno@usp.br\";\n\nlet verIndex = -1;\nlet copyrightIndex = -1;\n\nfor (let i = 0; i < items.length; i++) {\n  const text = items[i].text;\n  if (verIndex === -1 && text === VER_TEXT) {\n    verIndex = i;\n  }\n  if (copyrightIndex === -1 && text.indexOf(COPYRIGHT_MARKER) !== -1) {\n    copyrightIndex = i;\n  }\n}\n\nif (verIndex === -1 || copyrightIndex === -1) {\n  throw new Error(\n    \"edit.js: could not locate the 'Ver no Jupiter' / copyright paragraphs to remove\"\n  );\n}\n\n// Delete from the bottom up so earlier indices stay valid.\nitems[copyrightIndex].delete();\nitems[verIndex].delete();\n\n// The empty paragraph right before \"Ver no Jupiter ...\" is also removed by\n// the diff; only delete it if it is indeed blank so we don't eat content.\nconst blankIndex = verIndex - 1;\nif (blankIndex >= 0 && items[blankIndex].text === \"\") {\n  items[blankIndex].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter Salvar em pdf Salvar em docx\" line,\n# the \"\u00a9 2020 . Contact: ...\" credits line, and the blank paragraph that\n# separates them from the preceding \"Requisitos\" content, mirroring the\n# upstream Jekyll/Github Pages rebuild that dropped that boilerplate.\n\n$d = $word.ActiveDocument\n\n$verIndex = -1\n$copyrightIndex = -1\n$count = $d.Paragraphs.Count\n\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($verIndex -eq -1 -and $t -like \"Ver no Jupiter Salvar em pdf Salvar em docx*\") {\n        $verIndex = $i\n    }\n    if ($copyrightIndex -eq -1 -and $t -like \"*luizeleno@usp.br*\") {\n        $copyrightIndex = $i\n    }\n}\n\nif ($verIndex -eq -1 -or $copyrightIndex -eq -1) {\n    throw \"edit.ps1: could not locate the 'Ver no Jupiter' / copyright paragraphs to remove\"\n}\n\n# Delete bottom-up so the other, still-to-process index stays valid.\n$d.Paragraphs.Item($copyrightIndex).Range.Delete()\n$d.Paragraphs.Item($verIndex).Range.Delete()\n\n# The empty paragraph right before \"Ver no Jupiter ...\" is also removed by\n# the diff; only delete it if it is indeed blank so we don't eat content.\n$blankIndex = $verIndex - 1\nif ($blankIndex -ge 1) {\n    $blankText = $d.Paragraphs.Item($blankIndex).Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($blankText -eq \"\") {\n        $d.Paragraphs.Item($blankIndex).Range.Delete()\n    }\n}\n"}
